# edit.ps1 - apply the "Add files via upload" change to CMP73010 (1).docx
#
# Summary of the target edit:
#   1. The paragraph "Baz chan"+"ges" (split across two runs around a
#      hidden _GoBack bookmark) becomes a single run "Baz changes" -
#      the visible text is unchanged, only the run/bookmark structure is
#      cleaned up (the stray bookmark is removed from here).
#   2. The first trailing empty paragraph is kept as an empty paragraph.
#   3. The second trailing empty paragraph becomes:
#        "Batman says "Get to the chopper!""
#   4. A new paragraph is appended after that:
#        "Rory internally questions Batman's decision to impersonate
#         Arnold Schwarzenegger for a 7 second Instagram story."
#      built out of three runs (matching the original author's split),
#      with the _GoBack bookmark re-anchored to the very end of it.
#
# We build exact OOXML fragments and push them through Range.InsertXML
# so that run boundaries / the bookmark placement match precisely
# (plain Range.Text / InsertAfter calls get silently coalesced into a
# single run by the engine's run-normalisation on save).

$d = $word.ActiveDocument

$apos  = [char]0x27      # '  (use to embed literal single quotes in XML attrs)
$ldq   = [char]0x201C    # “
$rdq   = [char]0x201D    # ”
$rsq   = [char]0x2019    # '

$wns = "xmlns:w=" + $apos + "http://schemas.openxmlformats.org/wordprocessingml/2006/main" + $apos

# --- 1. "Baz chan" / "ges" -> single run "Baz changes", bookmark dropped ---
$pBaz = $d.Paragraphs.Item(5)
$xmlBaz = "<w:p " + $wns + "><w:r><w:t>Baz changes</w:t></w:r></w:p>"
$pBaz.Range.InsertXML($xmlBaz) | Out-Null

# --- 2. Paragraph 6 is already empty; leave it as-is (stays <w:p/>) ---

# --- 3. Paragraph 7 (currently empty) becomes the Batman line ---
$pBatman = $d.Paragraphs.Item(7)
$xmlBatman = "<w:p " + $wns + "><w:r><w:t>Batman says " + $ldq + "Get to the chopper!" + $rdq + "</w:t></w:r></w:p>"
$pBatman.Range.InsertXML($xmlBatman) | Out-Null

# --- 4. Append the new "Rory internally..." paragraph after it ---
$pBatman = $d.Paragraphs.Item(7)
$pBatman.Range.InsertParagraphAfter() | Out-Null
$pRory = $d.Paragraphs.Item(8)

$xmlRory = "<w:p " + $wns + ">" + `
    "<w:r><w:t>Rory internally questions Batman" + $rsq + "s decision to impersonate Arnold S</w:t></w:r>" + `
    "<w:r><w:t>chwarzenegger</w:t></w:r>" + `
    "<w:r><w:t xml:space=" + $apos + "preserve" + $apos + "> for a 7 second Instagram story.</w:t></w:r>" + `
    "<w:bookmarkStart w:id=" + $apos + "0" + $apos + " w:name=" + $apos + "_GoBack" + $apos + "/>" + `
    "<w:bookmarkEnd w:id=" + $apos + "0" + $apos + "/>" + `
    "</w:p>"

$pRory.Range.InsertXML($xmlRory) | Out-Null
